$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object 'object[,]' 24,5
$bf[0,0] = 1.02
$bf[0,1] = 1.025130049555415
$bf[0,2] = 1.052578524560983
$bf[0,3] = 1.025529291470578
$bf[0,4] = 1.055234760281801
$bf[1,0] = 1.02
$bf[1,1] = 1.025932597627786
$bf[1,2] = 1.05327366321171
$bf[1,3] = 1.026206640501578
$bf[1,4] = 1.056092013927693
$bf[2,0] = 1.02
$bf[2,1] = 1.026452506313345
$bf[2,2] = 1.053723630134669
$bf[2,3] = 1.02664584849754
$bf[2,4] = 1.056647431705833
$bf[3,0] = 1.02
$bf[3,1] = 1.026671219447121
$bf[3,2] = 1.053912833993863
$bf[3,3] = 1.026830709474855
$bf[3,4] = 1.056881099073213
$bf[4,0] = 1.02
$bf[4,1] = 1.0267079507128
$bf[4,2] = 1.053944604313257
$bf[4,3] = 1.026861761173653
$bf[4,4] = 1.056920342747442
$bf[5,0] = 1.02
$bf[5,1] = 1.026455428207069
$bf[5,2] = 1.053726158139913
$bf[5,3] = 1.026648317764903
$bf[5,4] = 1.05665055331621
$bf[6,0] = 1.02
$bf[6,1] = 1.025401147875437
$bf[6,2] = 1.05281341445761
$bf[6,3] = 1.025758013608556
$bf[6,4] = 1.055524323431224
$bf[7,0] = 1.02
$bf[7,1] = 1.023548088902783
$bf[7,2] = 1.051206397229244
$bf[7,3] = 1.024196293603702
$bf[7,4] = 1.0535453449062
$bf[8,0] = 1.02
$bf[8,1] = 1.022315986717273
$bf[8,2] = 1.050136076382554
$bf[8,3] = 1.023160034829829
$bf[8,4] = 1.05222989827848
$bf[9,0] = 1.02
$bf[9,1] = 1.021783269310765
$bf[9,2] = 1.049672882496396
$bf[9,3] = 1.022712505476886
$bf[9,4] = 1.051661238930626
$bf[10,0] = 1.02
$bf[10,1] = 1.02158551484574
$bf[10,2] = 1.049500872973274
$bf[10,3] = 1.022546451914051
$bf[10,4] = 1.051450156433925
$bf[11,0] = 1.02
$bf[11,1] = 1.021627928370202
$bf[11,2] = 1.049537767699537
$bf[11,3] = 1.022582062847735
$bf[11,4] = 1.051495427846809
$bf[12,0] = 1.02
$bf[12,1] = 1.021766920397219
$bf[12,2] = 1.049658663281684
$bf[12,3] = 1.022698775780119
$bf[12,4] = 1.051643787858698
$bf[13,0] = 1.02
$bf[13,1] = 1.021852573989493
$bf[13,2] = 1.049733156571602
$bf[13,3] = 1.022770710213274
$bf[13,4] = 1.051735216355249
$bf[14,0] = 1.02
$bf[14,1] = 1.022351358223964
$bf[14,2] = 1.050166822722466
$bf[14,3] = 1.0231897608467
$bf[14,4] = 1.052267658255958
$bf[15,0] = 1.02
$bf[15,1] = 1.022664445336285
$bf[15,2] = 1.050438921578749
$bf[15,3] = 1.023452936807043
$bf[15,4] = 1.052601897476375
$bf[16,0] = 1.02
$bf[16,1] = 1.022847139957333
$bf[16,2] = 1.050597657397856
$bf[16,3] = 1.023606556382985
$bf[16,4] = 1.052796943899829
$bf[17,0] = 1.02
$bf[17,1] = 1.022909446977007
$bf[17,2] = 1.050651786418795
$bf[17,3] = 1.023658955882064
$bf[17,4] = 1.052863464951643
$bf[18,0] = 1.02
$bf[18,1] = 1.02263084615786
$bf[18,2] = 1.050409725330246
$bf[18,3] = 1.023424688764903
$bf[18,4] = 1.052566027412232
$bf[19,0] = 1.02
$bf[19,1] = 1.021725987369807
$bf[19,2] = 1.049623061390129
$bf[19,3] = 1.022664401780862
$bf[19,4] = 1.051600095571869
$bf[20,0] = 1.02
$bf[20,1] = 1.021157764731544
$bf[20,2] = 1.049128694509056
$bf[20,3] = 1.022187414373168
$bf[20,4] = 1.050993603688672
$bf[21,0] = 1.02
$bf[21,1] = 1.02145892358492
$bf[21,2] = 1.049390744356596
$bf[21,3] = 1.022440175680241
$bf[21,4] = 1.051315037370012
$bf[22,0] = 1.02
$bf[22,1] = 1.022646027945071
$bf[22,2] = 1.050422917780707
$bf[22,3] = 1.023437452488799
$bf[22,4] = 1.052582235271318
$bf[23,0] = 1.02
$bf[23,1] = 1.024026579846842
$bf[23,2] = 1.051621677541154
$bf[23,3] = 1.024599182011914
$bf[23,4] = 1.054056284469986
$ws.Range("B2:F25").Value = $bf

$inl = New-Object 'object[,]' 24,6
$inl[0,0] = 1.04136903293982
$inl[0,1] = 1.030300792757304
$inl[0,2] = 1.055326868634535
$inl[0,3] = 1.028354595625322
$inl[0,4] = 1.057975777477524
$inl[0,5] = 1.014142848672106
$inl[1,0] = 1.04158885314406
$inl[1,1] = 1.030743301946246
$inl[1,2] = 1.055835183092671
$inl[1,3] = 1.028839619795867
$inl[1,4] = 1.058646320258604
$inl[1,5] = 1.014290135457412
$inl[2,0] = 1.04172991804204
$inl[2,1] = 1.03102958041562
$inl[2,2] = 1.056163584680788
$inl[2,3] = 1.029153701349418
$inl[2,4] = 1.0590802786068
$inl[2,5] = 1.014385393431295
$inl[3,0] = 1.04178893991622
$inl[3,1] = 1.031149917614937
$inl[3,2] = 1.056301520649422
$inl[3,3] = 1.029285796995087
$inl[3,4] = 1.059262730279219
$inl[3,5] = 1.01442542833231
$inl[4,0] = 1.0417988333977
$inl[4,1] = 1.031170121888378
$inl[4,2] = 1.056324673386414
$inl[4,3] = 1.029307979666292
$inl[4,4] = 1.05929336558121
$inl[4,5] = 1.014432149678271
$inl[5,0] = 1.041730707803626
$inl[5,1] = 1.031031188424749
$inl[5,2] = 1.056165428277329
$inl[5,3] = 1.029155466200476
$inl[5,4] = 1.059082716475925
$inl[5,5] = 1.01438592842554
$inl[6,0] = 1.041443564654985
$inl[6,1] = 1.030450351382793
$inl[6,2] = 1.055498760955622
$inl[6,3] = 1.02851846163445
$inl[6,4] = 1.058202374531494
$inl[6,5] = 1.014192634280827
$inl[7,0] = 1.040928633105286
$inl[7,1] = 1.029426477689088
$inl[7,2] = 1.054320154864144
$inl[7,3] = 1.027397863570214
$inl[7,4] = 1.056651730094896
$inl[7,5] = 1.013851689042731
$inl[8,0] = 1.040579381012677
$inl[8,1] = 1.028743721886721
$inl[8,2] = 1.053531917240485
$inl[8,3] = 1.02665214522112
$inl[8,4] = 1.055618494211709
$inl[8,5] = 1.013624192154155
$inl[9,0] = 1.040426747278912
$inl[9,1] = 1.02844805510218
$inl[9,2] = 1.053190027645584
$inl[9,3] = 1.026329577203249
$inl[9,4] = 1.055171238397338
$inl[9,5] = 1.013525641523467
$inl[10,0] = 1.040369842049699
$inl[10,1] = 1.028338228149476
$inl[10,2] = 1.053062949384766
$inl[10,3] = 1.026209812454791
$inl[10,4] = 1.055005130772822
$inl[10,5] = 1.013489029415551
$inl[11,0] = 1.0403820579156
$inl[11,1] = 1.028361786526436
$inl[11,2] = 1.053090211936098
$inl[11,3] = 1.026235500053306
$inl[11,4] = 1.055040760348881
$inl[11,5] = 1.01349688310213
$inl[12,0] = 1.040422047759258
$inl[12,1] = 1.028438976829129
$inl[12,2] = 1.053179525047037
$inl[12,3] = 1.026319676355908
$inl[12,4] = 1.055157507407984
$inl[12,5] = 1.013522615275376
$inl[13,0] = 1.040446658985478
$inl[13,1] = 1.028486535941892
$inl[13,2] = 1.053234542552377
$inl[13,3] = 1.026371547002361
$inl[13,4] = 1.055229442237222
$inl[13,5] = 1.013538468932624
$inl[14,0] = 1.040589481276967
$inl[14,1] = 1.028763343809842
$inl[14,2] = 1.053554595289423
$inl[14,3] = 1.026673560128567
$inl[14,4] = 1.055648180250743
$inl[14,5] = 1.013630731767135
$inl[15,0] = 1.040678694303853
$inl[15,1] = 1.028936971161658
$inl[15,2] = 1.053755202375723
$inl[15,3] = 1.026863095079635
$inl[15,4] = 1.055910882830404
$inl[15,5] = 1.013688594611117
$inl[16,0] = 1.040730595106768
$inl[16,1] = 1.029038242177308
$inl[16,2] = 1.053872157272655
$inl[16,3] = 1.02697367964389
$inl[16,4] = 1.056064126398235
$inl[16,5] = 1.013722340852876
$inl[17,0] = 1.040748268905149
$inl[17,1] = 1.029072772450901
$inl[17,2] = 1.053912026355383
$inl[17,3] = 1.027011391514838
$inl[17,4] = 1.056116380748032
$inl[17,5] = 1.013733846727457
$inl[18,0] = 1.040669136614683
$inl[18,1] = 1.028918342879239
$inl[18,2] = 1.053733684896894
$inl[18,3] = 1.026842756449931
$inl[18,4] = 1.055882695910734
$inl[18,5] = 1.013682386904344
$inl[19,0] = 1.040410277544425
$inl[19,1] = 1.028416246285949
$inl[19,2] = 1.053153226902845
$inl[19,3] = 1.026294887103744
$inl[19,4] = 1.055123127657277
$inl[19,5] = 1.013515037954236
$inl[20,0] = 1.040246306474652
$inl[20,1] = 1.028100540638208
$inl[20,2] = 1.052787777487921
$inl[20,3] = 1.025950717814004
$inl[20,4] = 1.054645691455366
$inl[20,5] = 1.013409784440112
$inl[21,0] = 1.040333345625237
$inl[21,1] = 1.028267903437353
$inl[21,2] = 1.052981555307096
$inl[21,3] = 1.026133139776
$inl[21,4] = 1.054898776087729
$inl[21,5] = 1.013465584450743
$inl[22,0] = 1.040673455742237
$inl[22,1] = 1.028926760207272
$inl[22,2] = 1.053743407892032
$inl[22,3] = 1.026851946502422
$inl[22,4] = 1.055895432324637
$inl[22,5] = 1.013685191912786
$inl[23,0] = 1.041062810340095
$inl[23,1] = 1.02969120938791
$inl[23,2] = 1.054625300084278
$inl[23,3] = 1.027687333415926
$inl[23,4] = 1.05705252349645
$inl[23,5] = 1.013939868724707
$ws.Range("I2:N25").Value = $inl
